$wb = $excel.ActiveWorkbook

# --- Sheet 1: DPLKKPS139-001 ---
$ws1 = $wb.Worksheets.Item("DPLKKPS139-001")

# Update the NO_REGISTER value (column N) from the old register number to the new one
$ws1.Range("N2").Value = "M03220800000029"

# Update the PREPARATION text block (column F) to reference the new register number
$ws1.Range("F2").Value = "Username : 31816;`nPassword : bni1234;`nRole : 09 - Penyelia Settlement;`nNo Register : M03220800000029;`nStatus Verifikasi : 0 : Kembalikan ke Register;`nKeterangan : KEP-012"

# Update the selected cell on this sheet's view to G2
$ws1.Activate()
[void]$ws1.Range("G2").Select()

# --- Sheet 2: DPLKKPS139-002 ---
$ws2 = $wb.Worksheets.Item("DPLKKPS139-002")

# Update the selected cell on this sheet's view back to A2 and make it the active sheet
$ws2.Activate()
[void]$ws2.Range("A2").Select()
